$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Issue #57: Make genre required, with PBCore controlled vocabulary.
# The fixture spreadsheet gains a new "Genre" column (T) to exercise the
# manifest-driven genre override: header in row 2, and one controlled
# vocabulary term per existing data row (rows 3-4).

$ws.Range("T2").Value = "Genre"
$ws.Range("T3").Value = "Travel"
$ws.Range("T4").Value = "Bicycle"

# Match the refreshed selection left after editing the sheet.
$null = $ws.Range("T4").Select()
